$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.449.05"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.163.44"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.23"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.20"
$ws.Range("E7").Value = "  +4.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.396"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0863"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.01"
$ws.Range("E12").Value = "  +6.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.483.04"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.20"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.816"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.58"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.166.57"
$ws.Range("E17").Value = "  +3.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.434.54"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.25"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.95"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.66"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.00"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.73"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  +9.90%  "
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("E35").Value = "  +9.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0624"
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("E37").Value = "  +2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.61"
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "104.61"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.00"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.538.56"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +6.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0935"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.91"
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("E47").Value = "  +7.42%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.21"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.367.69"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("E51").Value = "  +0.37%  "
